$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap F:V content between row 66 and row 67
$ws.Range("F66").Value = 'Chaves'
$ws.Range("G66").Value = 4
$ws.Range("H66").Value = 'Gil Vicente'
$ws.Range("I66").Value = 2
$ws.Range("J66").Value = 2.62
$ws.Range("K66").Value = '02/10/2023 20:42'
$ws.Range("L66").Value = 2.74
$ws.Range("M66").Value = '07/10/2023 16:29'
$ws.Range("N66").Value = 3.45
$ws.Range("O66").Value = '02/10/2023 20:42'
$ws.Range("P66").Value = 3.65
$ws.Range("Q66").Value = '07/10/2023 16:25'
$ws.Range("R66").Value = 2.73
$ws.Range("S66").Value = '02/10/2023 20:42'
$ws.Range("T66").Value = 2.58
$ws.Range("U66").Value = '07/10/2023 16:29'
$ws.Range("V66").Value = 'https://www.betexplorer.com/football/portugal/liga-portugal/chaves-gil-vicente/K4BKKZh1/'
$ws.Range("F67").Value = 'SC Farense'
$ws.Range("G67").Value = 0
$ws.Range("H67").Value = 'Vizela'
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 2.29
$ws.Range("K67").Value = '02/10/2023 07:12'
$ws.Range("L67").Value = 2.32
$ws.Range("M67").Value = '07/10/2023 16:02'
$ws.Range("N67").Value = 3.4
$ws.Range("O67").Value = '02/10/2023 07:12'
$ws.Range("P67").Value = 3.47
$ws.Range("Q67").Value = '07/10/2023 15:49'
$ws.Range("R67").Value = 3.31
$ws.Range("S67").Value = '02/10/2023 07:12'
$ws.Range("T67").Value = 3.24
$ws.Range("U67").Value = '07/10/2023 16:03'
$ws.Range("V67").Value = 'https://www.betexplorer.com/football/portugal/liga-portugal/sc-farense-vizela/OY1Asc0E/'

# Swap F:V content between row 76 and row 77
$ws.Range("F76").Value = 'Benfica'
$ws.Range("G76").Value = 1
$ws.Range("H76").Value = 'Casa Pia'
$ws.Range("I76").Value = 1
$ws.Range("J76").Value = 1.22
$ws.Range("K76").Value = '11/10/2023 14:42'
$ws.Range("L76").Value = 1.22
$ws.Range("M76").Value = '28/10/2023 18:55'
$ws.Range("N76").Value = 7.32
$ws.Range("O76").Value = '11/10/2023 14:42'
$ws.Range("P76").Value = 6.95
$ws.Range("Q76").Value = '28/10/2023 18:58'
$ws.Range("R76").Value = 13.22
$ws.Range("S76").Value = '11/10/2023 14:42'
$ws.Range("T76").Value = 14.17
$ws.Range("U76").Value = '28/10/2023 18:58'
$ws.Range("V76").Value = 'https://www.betexplorer.com/football/portugal/liga-portugal/benfica-casa-pia/GWtkzFhl/'
$ws.Range("F77").Value = 'Vitoria Guimaraes'
$ws.Range("G77").Value = 5
$ws.Range("H77").Value = 'Chaves'
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 1.62
$ws.Range("K77").Value = '11/10/2023 14:42'
$ws.Range("L77").Value = 1.69
$ws.Range("M77").Value = '28/10/2023 18:58'
$ws.Range("N77").Value = 4.26
$ws.Range("O77").Value = '11/10/2023 14:42'
$ws.Range("P77").Value = 4
$ws.Range("Q77").Value = '28/10/2023 18:58'
$ws.Range("R77").Value = 5.63
$ws.Range("S77").Value = '11/10/2023 14:42'
$ws.Range("T77").Value = 5.29
$ws.Range("U77").Value = '28/10/2023 18:58'
$ws.Range("V77").Value = 'https://www.betexplorer.com/football/portugal/liga-portugal/vitoria-guimaraes-chaves/8vH9wlat/'

# Swap F:V content between row 87 and row 88
$ws.Range("F87").Value = 'Casa Pia'
$ws.Range("G87").Value = 0
$ws.Range("H87").Value = 'Vizela'
$ws.Range("I87").Value = 1
$ws.Range("J87").Value = 2.17
$ws.Range("K87").Value = '29/10/2023 21:42'
$ws.Range("L87").Value = 2.25
$ws.Range("M87").Value = '05/11/2023 16:22'
$ws.Range("N87").Value = 3.42
$ws.Range("O87").Value = '29/10/2023 21:42'
$ws.Range("P87").Value = 3.28
$ws.Range("Q87").Value = '05/11/2023 16:21'
$ws.Range("R87").Value = 3.58
$ws.Range("S87").Value = '29/10/2023 21:42'
$ws.Range("T87").Value = 3.59
$ws.Range("U87").Value = '05/11/2023 16:22'
$ws.Range("V87").Value = 'https://www.betexplorer.com/football/portugal/liga-portugal/casa-pia-vizela/6gDUYPDB/'
$ws.Range("F88").Value = 'Moreirense'
$ws.Range("G88").Value = 1
$ws.Range("H88").Value = 'Vitoria Guimaraes'
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 2.66
$ws.Range("K88").Value = '30/10/2023 14:42'
$ws.Range("L88").Value = 2.93
$ws.Range("M88").Value = '05/11/2023 16:28'
$ws.Range("N88").Value = 3.23
$ws.Range("O88").Value = '30/10/2023 14:42'
$ws.Range("P88").Value = 3.35
$ws.Range("Q88").Value = '05/11/2023 16:02'
$ws.Range("R88").Value = 2.9
$ws.Range("S88").Value = '30/10/2023 14:42'
$ws.Range("T88").Value = 2.58
$ws.Range("U88").Value = '05/11/2023 16:21'
$ws.Range("V88").Value = 'https://www.betexplorer.com/football/portugal/liga-portugal/moreirense-vitoria-guimaraes/vq5sWNcU/'

# Add new row 90 - copy formatting from row 89 first
$ws.Range("A89:V89").Copy()
$ws.Range("A90:V90").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A90").Value = 89
$ws.Range("B90").Value = 'portugal'
$ws.Range("C90").Value = 'liga-portugal'
$ws.Range("D90").Value = '2023-2024'
$ws.Range("E90").Value = 45235.89583333334
$ws.Range("F90").Value = 'Sporting CP'
$ws.Range("G90").Value = 3
$ws.Range("H90").Value = 'Estrela'
$ws.Range("I90").Value = 2
$ws.Range("J90").Value = 1.2
$ws.Range("K90").Value = '30/10/2023 21:42'
$ws.Range("L90").Value = 1.19
$ws.Range("M90").Value = '05/11/2023 21:23'
$ws.Range("N90").Value = 7.62
$ws.Range("O90").Value = '30/10/2023 21:42'
$ws.Range("P90").Value = 7.78
$ws.Range("Q90").Value = '05/11/2023 21:23'
$ws.Range("R90").Value = 12.34
$ws.Range("S90").Value = '30/10/2023 21:42'
$ws.Range("T90").Value = 14.32
$ws.Range("U90").Value = '05/11/2023 21:23'
$ws.Range("V90").Value = 'https://www.betexplorer.com/football/portugal/liga-portugal/sporting-lisbon-estrela-da-amadora/CWZekobh/'
